# Append the sheet's id number, e.g. "almere" -> "almere id(1)"
$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $sheet = $wb.Worksheets.Item($i)
    $sheet.Name = "$($sheet.Name) id($i)"
}
